$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $style = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $style
}

# Row 2
Set-TextValue $ws.Range("D2") "30.109.56"
Set-TextValue $ws.Range("E2") "  +4.40%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.908.56"
Set-TextValue $ws.Range("E3") "  +5.37%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.000"
Set-TextValue $ws.Range("E4") "  +0.11%  "

# Row 5
Set-TextValue $ws.Range("D5") "251.21"
Set-TextValue $ws.Range("E5") "  +0.79%  "

# Row 6
Set-TextValue $ws.Range("D6") "0.9999"
Set-TextValue $ws.Range("E6") "  +0.03%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.5113"
Set-TextValue $ws.Range("E7") "  +3.51%  "

# Row 8
Set-TextValue $ws.Range("D8") "44.92"
Set-TextValue $ws.Range("E8") "  +3.80%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.2965"
Set-TextValue $ws.Range("E9") "  +6.44%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.06807"
Set-TextValue $ws.Range("E10") "  +5.58%  "

# Row 11
Set-TextValue $ws.Range("D11") "1.910.95"
Set-TextValue $ws.Range("E11") "  +5.50%  "

# Row 12
Set-TextValue $ws.Range("D12") "17.24"
Set-TextValue $ws.Range("E12") "  +2.47%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.07362"
Set-TextValue $ws.Range("E13") "  +3.59%  "

# Row 14
Set-TextValue $ws.Range("D14") "0.6908"
Set-TextValue $ws.Range("E14") "  +6.46%  "

# Row 15
Set-TextValue $ws.Range("D15") "86.74"
Set-TextValue $ws.Range("E15") "  +2.60%  "

# Row 16
Set-TextValue $ws.Range("D16") "4.876"
Set-TextValue $ws.Range("E16") "  +3.50%  "

# Row 17
Set-TextValue $ws.Range("B17") "WrappedBTC"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D17") "30.119.90"
Set-TextValue $ws.Range("E17") "  +4.50%  "

# Row 18
Set-TextValue $ws.Range("B18") "ShibaInu"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D18") "0.000008294"
Set-TextValue $ws.Range("E18") "  +12.13%  "

# Row 19
Set-TextValue $ws.Range("D19") "1.000"
Set-TextValue $ws.Range("E19") "  +0.05%  "

# Row 20
Set-TextValue $ws.Range("D20") "12.96"
Set-TextValue $ws.Range("E20") "  +5.81%  "

# Row 21
Set-TextValue $ws.Range("D21") "2.156.26"
Set-TextValue $ws.Range("E21") "  +5.27%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.9999"
Set-TextValue $ws.Range("E22") "  +0.08%  "

# Row 23
Set-TextValue $ws.Range("D23") "4.812"
Set-TextValue $ws.Range("E23") "  +4.65%  "

# Row 24
Set-TextValue $ws.Range("D24") "5.705"
Set-TextValue $ws.Range("E24") "  +6.53%  "

# Row 25
Set-TextValue $ws.Range("D25") "9.139"
Set-TextValue $ws.Range("E25") "  +2.20%  "

# Row 26
Set-TextValue $ws.Range("D26") "146.82"
Set-TextValue $ws.Range("E26") "  +2.11%  "

# Row 27
Set-TextValue $ws.Range("D27") "134.92"
Set-TextValue $ws.Range("E27") "  +1.28%  "

# Row 28
Set-TextValue $ws.Range("D28") "16.99"
Set-TextValue $ws.Range("E28") "  +1.78%  "

# Row 29
Set-TextValue $ws.Range("D29") "1.999"
Set-TextValue $ws.Range("E29") "  +5.49%  "

# Row 30
Set-TextValue $ws.Range("E30") "  +0.06%  "

# Row 31
Set-TextValue $ws.Range("D31") "4.225"
Set-TextValue $ws.Range("E31") "  +1.43%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.08812"
Set-TextValue $ws.Range("E32") "  +5.44%  "

# Row 33
Set-TextValue $ws.Range("D33") "4.006"
Set-TextValue $ws.Range("E33") "  +4.42%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.05061"
Set-TextValue $ws.Range("E34") "  +1.57%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.146"
Set-TextValue $ws.Range("E35") "  +4.85%  "

# Row 36
Set-TextValue $ws.Range("E36") "  +5.61%  "

# Row 37
Set-TextValue $ws.Range("D37") "2.688"
Set-TextValue $ws.Range("E37") "  -0.27%  "

# Row 38
Set-TextValue $ws.Range("D38") "2.809"
Set-TextValue $ws.Range("E38") "  +1.96%  "

# Row 39
Set-TextValue $ws.Range("D39") "2.276"
Set-TextValue $ws.Range("E39") "  +0.14%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.9689"
Set-TextValue $ws.Range("E40") "  +1.20%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.01691"
Set-TextValue $ws.Range("E41") "  +6.05%  "

# Row 42
Set-TextValue $ws.Range("D42") "6.158"
Set-TextValue $ws.Range("E42") "  +1.77%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.4292"
Set-TextValue $ws.Range("E43") "  +4.72%  "

# Row 44
Set-TextValue $ws.Range("D44") "105.06"
Set-TextValue $ws.Range("E44") "  +4.93%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.9989"
Set-TextValue $ws.Range("E45") "  -0.06%  "

# Row 46
Set-TextValue $ws.Range("D46") "7.617"
Set-TextValue $ws.Range("E46") "  +5.46%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.1276"
Set-TextValue $ws.Range("E47") "  +4.16%  "

# Row 48
Set-TextValue $ws.Range("D48") "0.05733"
Set-TextValue $ws.Range("E48") "  +3.62%  "

# Row 49
Set-TextValue $ws.Range("E49") "  +4.56%  "

# Row 50
Set-TextValue $ws.Range("D50") "8.410"
Set-TextValue $ws.Range("E50") "  +3.64%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.3796"
Set-TextValue $ws.Range("E51") "  +4.50%  "
